# Modif mineur au excel
#
# 1) Highlight (yellow fill) the header rows for the two "Marx" / "Romane"
#    groups (rows 15 & 37, columns A:D) on sheet "Individus".
# 2) Remove three unused, empty columns (Q:S) from sheet "Cheveux" so the
#    summary table that used to start at column U now starts at column R.
# 3) Restore the selections / active cell that Excel leaves behind after
#    those edits.

$wb = $excel.ActiveWorkbook

$wsIndividus = $wb.Worksheets.Item("Individus")
$wsCheveux   = $wb.Worksheets.Item("Cheveux")

# --- 1) Yellow-highlight the two header rows -----------------------------
$wsIndividus.Range("A15:D15").Interior.Color = 65535
$wsIndividus.Range("A37:D37").Interior.Color = 65535

# --- 2) Delete the 3 empty columns (Q:S) on "Cheveux" ---------------------
# Columns Q, R, S are unused (data resumes at U); deleting them shifts the
# U..AF summary block left by 3 columns, landing on R..AC.
$wsCheveux.Range("Q1:S1").EntireColumn.Delete()

# --- 3) Restore selections -------------------------------------------------
$wsCheveux.Select()
$wsCheveux.Range("W25").Select()

$wsIndividus.Select()
$wsIndividus.Range("E28").Select()
